$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 453; this shifts the existing rows 453-503
# down to 454-504 (exactly matching the diff's row-shift pattern).
$ws.Rows("453:453").Insert()

# Populate the newly inserted row 453 with a new data record. Most fields
# mirror the constant values shared by every row in this data block; only
# D (Fecha), K/L/M (prices) and P (Precio $/Kg) carry new figures.
$ws.Range("A453").Value = 5
$ws.Range("B453").Value = "Macroferia Regional de Talca"
$ws.Range("C453").Value = "Maule"
$ws.Range("D453").Value = 45194
$ws.Range("E453").Value = 7
$ws.Range("F453").Value = 100112009
$ws.Range("G453").Value = "Acelga"
$ws.Range("H453").Value = "Sin especificar"
$ws.Range("I453").Value = "Primera"
$ws.Range("J453").Value = 500
$ws.Range("K453").Value = 2000
$ws.Range("L453").Value = 2000
$ws.Range("M453").Value = 2000
$ws.Range("N453").Value = "$/docena de atados (4 kilos)"
$ws.Range("O453").Value = "Región del Maule"
$ws.Range("P453").Value = 500
$ws.Range("Q453").Value = 4
$ws.Range("R453").Value = "Hortaliza"
